$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 <-> Row 3 swap of D (Fecha) and J (Volumen)
$ws.Range("D2").Value = 44277
$ws.Range("J2").Value = 150

$ws.Range("D3").Value = 44291
$ws.Range("J3").Value = 30

# Row 4 <-> Row 5 swap of D (Fecha)
$ws.Range("D4").Value = 44284
$ws.Range("D5").Value = 44280
